# Updated cryptos list on Mon Sep 30 03:32:35 UTC 2024 with GitHub Actions
#
# Price cells (column D) hold text-formatted numbers (e.g. thousand-dot
# separators like "64.418.04"). Plain numeric-looking strings such as
# "577.89" would otherwise be auto-detected by Excel as a Number when
# assigned via .Value, so those are written with a leading apostrophe to
# force text storage, matching the source data's text cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.418.04"
$ws.Range("E2").Value = "  -2.15%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.621.63"
$ws.Range("E3").Value = "  -1.96%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'577.89"
$ws.Range("E5").Value = "  -3.76%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'156.63"
$ws.Range("E6").Value = "  -0.58%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.645"
$ws.Range("E7").Value = "  +5.56%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -5.62%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "'5.76"
$ws.Range("E10").Value = "  -1.70%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -2.68%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.31%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'28.41"
$ws.Range("E13").Value = "  -2.58%  "

# Row 14 - ShibaInu
$ws.Range("D14").Value = "'0.0000186"
$ws.Range("E14").Value = "  -7.83%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.092.35"
$ws.Range("E15").Value = "  -1.92%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "64.285.57"
$ws.Range("E16").Value = "  -2.08%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.617.05"
$ws.Range("E17").Value = "  -2.36%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "'12.24"
$ws.Range("E18").Value = "  -4.05%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "'4.68"
$ws.Range("E19").Value = "  -2.96%  "

# Row 20 - now Uniswap (was BitcoinCash)
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'7.36"
$ws.Range("E20").Value = "  -2.35%  "

# Row 21 - now BitcoinCash (was Uniswap)
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'344.81"
$ws.Range("E21").Value = "  -2.25%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.25%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'67.91"
$ws.Range("E23").Value = "  -2.46%  "

# Row 24 - PEPE
$ws.Range("D24").Value = "'0.0000110"
$ws.Range("E24").Value = "  -6.06%  "

# Row 25 - SuiNetwork
$ws.Range("D25").Value = "'1.72"
$ws.Range("E25").Value = "  +2.48%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "'9.38"
$ws.Range("E26").Value = "  -3.39%  "

# Row 27 - Fetch.AI
$ws.Range("D27").Value = "'1.55"
$ws.Range("E27").Value = "  -3.15%  "

# Row 28 - Bittensor
$ws.Range("D28").Value = "'553.56"
$ws.Range("E28").Value = "  +3.39%  "

# Row 29 - Kaspa
$ws.Range("E29").Value = "  -2.57%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  +0.19%  "

# Row 31 - Aptos
$ws.Range("D31").Value = "'7.91"
$ws.Range("E31").Value = "  -1.45%  "

# Row 32 - PancakeSwap
$ws.Range("D32").Value = "'2.07"
$ws.Range("E32").Value = "  -3.19%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  -2.33%  "

# Row 34 - RenderToken
$ws.Range("D34").Value = "'6.40"
$ws.Range("E34").Value = "  -1.68%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "'5.30"
$ws.Range("E35").Value = "  -3.28%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("E36").Value = "  -2.86%  "

# Row 37 - EthereumClassic
$ws.Range("D37").Value = "'20.00"
$ws.Range("E37").Value = "  -3.16%  "

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "  -0.04%  "

# Row 39 - Stacks
$ws.Range("D39").Value = "'1.93"
$ws.Range("E39").Value = "  -0.95%  "

# Row 40 - Monero
$ws.Range("D40").Value = "'151.69"
$ws.Range("E40").Value = "  -4.15%  "

# Row 41 - USDe
$ws.Range("E41").Value = "  -0.03%  "

# Row 42 - dogwifhat
$ws.Range("D42").Value = "'2.44"
$ws.Range("E42").Value = "  +2.71%  "

# Row 43 - Aave
$ws.Range("D43").Value = "'158.56"
$ws.Range("E43").Value = "  -3.07%  "

# Row 44 - Filecoin
$ws.Range("D44").Value = "'3.98"
$ws.Range("E44").Value = "  -3.33%  "

# Row 45 - Hedera
$ws.Range("D45").Value = "'0.0603"
$ws.Range("E45").Value = "  -1.63%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = "'22.74"
$ws.Range("E46").Value = "  -0.06%  "

# Row 47 - Mantle
$ws.Range("D47").Value = "'0.632"
$ws.Range("E47").Value = "  -1.35%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  +2.42%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  -3.18%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "'19.09"
$ws.Range("E50").Value = "  -4.81%  "

# Row 51 - BabyDogeCoin
$ws.Range("D51").Value = "0.0₆0238"
$ws.Range("E51").Value = "  -8.10%  "
